$wb = $excel.ActiveWorkbook

$wsIndex = $wb.Worksheets.Item("index")
$wsP1 = $wb.Worksheets.Item("p1")

# --- Content edits ---
# New strings must be written in this order so the rebuilt shared-string
# table lands on the same indices as the target file:
#   1) index!B5 link target (draft:guide:top2 -> mdl:guide41:top)
#   2) p1!B7 placeholder text ("work in progress" notice)
#   3) p1!B2 page title suffix

$wsIndex.Range("B5").Value = '<a class="btn btn-primary btn-xs" role="button" href="https://support.vle.hiroshima-u.ac.jp/mdl:guide41:top" style="width:45%">マニュアルの目次へ戻る</a><br>'

$prepText = "現在作成中です。" + [char]10 + "順次公開いたしますので、公開前のページについては過去のマニュアルを参照してください。" + [char]10 + "■過去のマニュアル" + [char]10 + "https://support.vle.hiroshima-u.ac.jp/files/public/hirodai-moodle-faculty-document-20230306.pdf"
$wsP1.Range("B7").Value = $prepText
$wsP1.Rows.Item(7).RowHeight = 54

$wsP1.Range("B2").Value = "広大moodleへアクセスする【準備中】"

# --- View / selection edits ---
# Set the "index" sheet's remembered selection first (it is still the
# active sheet at this point), then switch to / select on "p1" last so it
# ends up as the active tab with its own remembered selection.
$wsIndex.Range("B17").Select() | Out-Null

$wsP1.Activate() | Out-Null
$wsP1.Range("B11").Select() | Out-Null
